# v1 final edits; creating tag
#
# XEM6310 sheet: the "IOStandard" column (W) used two different literal
# spellings for the same two standards ("LVCMOS" / "lvcmos" and "lvds").
# Clean that up by standardizing on the FPGA-style part names ("LVCMOS33",
# "LVDS_25"), label the column "IOStandard", and add a new "DiffTerm"
# column (W -> X) flagging the differential-pair data/clock-out lines that
# need termination enabled.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XEM6310")

# New header for the (previously unlabeled) IOStandard column.
$ws.Cells.Item(1, 23).Value = "IOStandard"

# Rows whose W column said "LVCMOS" or "lvcmos" -> "LVCMOS33".
$rowsLvcmos = @(
    12, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 28, 29, 30, 31, 32, 33, 34,
    38, 40, 42, 44, 46, 48, 50, 52, 58, 59, 60, 61, 89, 91, 93, 113, 115,
    138, 153, 155
)
foreach ($r in $rowsLvcmos) {
    $ws.Cells.Item($r, 23).Value = "LVCMOS33"
}

# Rows whose W column said "lvds" -> "LVDS_25".
$rowsLvds = @(
    104, 106, 108, 110, 112, 114, 118, 120, 127, 129, 131, 133, 135, 139,
    140, 141, 142, 143, 144, 145, 146, 147, 148, 149, 150, 151, 152, 154,
    156, 157, 158, 160
)
foreach ($r in $rowsLvds) {
    $ws.Cells.Item($r, 23).Value = "LVDS_25"
}

# New "DiffTerm" column header, and TRUE for the ADC data (_D_P/_D_N) and
# clock-out (_DCO_P/_DCO_N) differential pairs.
$ws.Cells.Item(1, 24).Value = "DiffTerm"
$rowsDiffTerm = @(108, 110, 112, 114, 131, 133, 135, 139, 145, 147, 148, 149, 150, 151, 152, 154)
foreach ($r in $rowsDiffTerm) {
    $ws.Cells.Item($r, 24).Value = $true
}

# Leave the selection where the author left it when saving.
$ws.Range("X13").Select() | Out-Null
